$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old blank placeholder row (row 3) so the data block shifts up
#    by one row. This matches the final layout where the sorted data starts
#    at row 3 (right under the header row 2) instead of row 4.
# ---------------------------------------------------------------------------
$ws.Rows(3).Delete()

# ---------------------------------------------------------------------------
# 2. Re-write the 9 sector rows (now rows 3-11) in descending order of income
#    (column D / E), and add the newly-computed analysis columns E:H.
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "种植业 "
$ws.Range("D3").Value = 1645.53
$ws.Range("E3").Value = 56.729111794475799
$ws.Range("F3").Value = 56.729111794475799
$ws.Range("G3").Value = 11.1111111111111
$ws.Range("H3").Value = 100

$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "畜牧业 "
$ws.Range("D4").Value = 390.24
$ws.Range("E4").Value = 13.453397134465
$ws.Range("F4").Value = 70.182508928940806
$ws.Range("G4").Value = 22.2222222222222
$ws.Range("H4").Value = 100

$ws.Range("B5").Value = 8
$ws.Range("C5").Value = "商饮服务业 "
$ws.Range("D5").Value = 211.62
$ws.Range("E5").Value = 7.2955307031454701
$ws.Range("F5").Value = 77.478039632086293
$ws.Range("G5").Value = 33.3333333333333
$ws.Range("H5").Value = 100

$ws.Range("B6").Value = 5
$ws.Range("C6").Value = "工业 "
$ws.Range("D6").Value = 167.38
$ws.Range("E6").Value = 5.7703710854006598
$ws.Range("F6").Value = 83.248410717486905
$ws.Range("G6").Value = 44.4444444444444
$ws.Range("H6").Value = 100

$ws.Range("B7").Value = 7
$ws.Range("C7").Value = "运输业 "
$ws.Range("D7").Value = 150.88
$ws.Range("E7").Value = 5.20153894948771
$ws.Range("F7").Value = 88.449949666974703
$ws.Range("G7").Value = 55.5555555555556
$ws.Range("H7").Value = 100

$ws.Range("B8").Value = 9
$ws.Range("C8").Value = "其他 "
$ws.Range("D8").Value = 136.69999999999999
$ws.Range("E8").Value = 4.7126880593516001
$ws.Range("F8").Value = 93.162637726326196
$ws.Range("G8").Value = 66.6666666666667
$ws.Range("H8").Value = 100

$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "林业 "
$ws.Range("D9").Value = 79.66
$ws.Range("E9").Value = 2.7462526028379601
$ws.Range("F9").Value = 95.908890329164194
$ws.Range("G9").Value = 77.7777777777778
$ws.Range("H9").Value = 100

$ws.Range("B10").Value = 4
$ws.Range("C10").Value = "渔业 "
$ws.Range("D10").Value = 74.12
$ws.Range("E10").Value = 2.5552629038708199
$ws.Range("F10").Value = 98.464153233035006
$ws.Range("G10").Value = 88.8888888888889
$ws.Range("H10").Value = 100

$ws.Range("B11").Value = 6
$ws.Range("C11").Value = "建筑业 "
$ws.Range("D11").Value = 44.55
$ws.Range("E11").Value = 1.53584676696499
$ws.Range("F11").Value = 100
$ws.Range("G11").Value = 100
$ws.Range("H11").Value = 100

# ---------------------------------------------------------------------------
# 3. Concentration-index helper columns (I:L), only populated on row 3.
# ---------------------------------------------------------------------------
$ws.Range("I3").Value = 763.62370202848899
$ws.Range("I3").NumberFormat = "0.00E+00"
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 900
$ws.Range("L3").Value = 0.65905925507122398

# ---------------------------------------------------------------------------
# 4. Hidden `_FilterDatabase` defined name (leftover from Data > Sort /
#    Filter usage) pointing at the data block.
# ---------------------------------------------------------------------------
$name = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$B`$4:`$E`$13")
$name.Visible = $false

# ---------------------------------------------------------------------------
# 5. Stamp a `sortState` record for the descending sort on column E that
#    produced the row order above. Row 12 (totals) is still blank at this
#    point, so including it in the sort range does not disturb the order
#    already written into rows 3-11.
# ---------------------------------------------------------------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("E2:E12"), 0, 2, 0, 0) | Out-Null
$ws.Sort.SetRange($ws.Range("B2:H12"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# ---------------------------------------------------------------------------
# 6. Totals row (row 12).
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = "合计"
$ws.Range("C12").Value = $null
$ws.Range("D13").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("D12").Value = 2900.68
$ws.Range("D12").NumberFormat = "0.00E+00"
$ws.Range("E12").Value = 100

# ---------------------------------------------------------------------------
# 7. New trailing blank row (row 13) with the same per-column styling as the
#    rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("B13").ClearContents()

$ws.Range("C3").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("C13").ClearContents()
$ws.Range("D13").ClearContents()

# ---------------------------------------------------------------------------
# 8. Cosmetic sheet-level tweaks.
# ---------------------------------------------------------------------------
$ws.Rows(2).RowHeight = 25.5
$ws.Range("H11").Select()
